$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add new "posWinPercent-ish" value and the new "Up" verdict string ---
$ws.Cells.Item(3, 24).Value = 0.069999999999993179
$ws.Cells.Item(3, 25).Value = "Up"

# --- Row 4: new trade entry ---
# Copy formatting (date style) from A3 into A4, then set the value
$ws.Cells.Item(3, 1).Copy($ws.Cells.Item(4, 1))
$ws.Cells.Item(4, 1).Value = 42633.884328703702

$ws.Cells.Item(4, 2).Value = -6
$ws.Cells.Item(4, 3).Value = "Neutral"
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = "Random"
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = -31.95

# Copy formatting (percentage style) from S3 into S4, then set the value
$ws.Cells.Item(3, 19).Copy($ws.Cells.Item(4, 19))
$ws.Cells.Item(4, 19).Value = -0.082000000000000003

$ws.Cells.Item(4, 20).Value = -0.28000000000000003
$ws.Cells.Item(4, 21).Value = 6.77
$ws.Cells.Item(4, 22).Value = 1.88
$ws.Cells.Item(4, 23).Value = 0
